$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for all existing data rows (2-422)
for ($r = 2; $r -le 422; $r++) {
    $ws.Cells.Item($r, 3).Value = 45182
}

# Row 422 gains an explicit row height (ht="15" customHeight="1")
$ws.Rows.Item(422).RowHeight = 15

# Append the new record as row 423
$row = 423
$ws.Cells.Item($row, 1).Value = "A 42391-2023"

$ws.Cells.Item($row, 2).Value = 45180
$ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($row, 3).Value = 45182
$ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($row, 4).Value = "DALARNAS LÄN"
$ws.Cells.Item($row, 5).Value = "MORA"
$ws.Cells.Item($row, 6).Value = "Kommuner"
$ws.Cells.Item($row, 7).Value = 0.3
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 0
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0

$ws.Cells.Item($row, 18).Value = ""
$ws.Cells.Item($row, 18).WrapText = $true
